$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)              # wdCollapseEnd -> collapse to end of title paragraph
$titleRange.InsertParagraphAfter()   # creates a brand-new empty paragraph after the title

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range

# Build the paragraph's contents explicitly via OOXML so the run layout
# matches exactly: an empty leading run, a bold "Meta description" run,
# and a plain run with the remaining text.
$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r/>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
  '<w:r><w:t>: Read our review of Cave of Fortune and play for free. With high volatility and big payouts, it' + [char]0x27 + 's perfect for experienced players.</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$null = $metaRange.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Remove the duplicate bold "Play Cave of Fortune Free..." paragraph that
#    sits right before the final (italic) paragraph.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupPara = $d.Paragraphs.Item($count - 1)
$dupPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new prompt,
#    keeping its existing run formatting (italic) untouched.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastRange = $lastPara.Range
# Exclude the trailing paragraph mark so only the visible text is replaced.
$textRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$textRange.Text = "Create an eye-catching feature image for Cave of Fortune slot game that features a happy Maya warrior wearing glasses. The image should be in a cartoon style that captures the game" + [char]0x27 + "s diamond mine setting and the Maya warrior" + [char]0x27 + "s joyous expression. Use bright and vivid colors that evoke the spirit of adventure and discovery. Give the warrior a fun, playful pose that invites players to join in on the excitement. Make sure the image accurately represents the game" + [char]0x27 + "s theme and captures its essence."
